# A new weekly price record was inserted into the "Ajo" (garlic) dataset.
# It belongs chronologically right after the existing row 89 (date 44195)
# and before the (old) row 90 (date 44447), so it is inserted as the new
# row 90 - pushing every subsequent record down by one row (old row 90
# becomes row 91, ..., old row 174 becomes row 175).
#
# The new record reuses the same market/category/variety/quality/price
# figures as the (old) row 90 record, only the date changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 90; everything below shifts down by one.
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new weekly record.
$ws.Range("A90").Value = 5
$ws.Range("B90").Value = "Macroferia Regional de Talca"
$ws.Range("C90").Value = "Maule"
$ws.Range("D90").Value = 44484
$ws.Range("E90").Value = 7
$ws.Range("F90").Value = 100112003
$ws.Range("G90").Value = "Ajo"
$ws.Range("H90").Value = "Chino"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 200
$ws.Range("K90").Value = 14000
$ws.Range("L90").Value = 14000
$ws.Range("M90").Value = 14000
$ws.Range("N90").Value = "`$/caja 10 kilos"
$ws.Range("O90").Value = "China"
$ws.Range("P90").Value = 1400
$ws.Range("Q90").Value = 10
$ws.Range("R90").Value = "Hortaliza"
